# Lesson 19 (3rd ed.) wordlist update:
# Insert 9 new "Useful Expressions: Greeting Cards" entries after row 52
# ("Please give my best regards..."), pushing the existing "spring" ...
# "decision by majority" block (old rows 53-110) down to rows 62-119.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new rows before row 53, shifting existing rows 53:110 down to 62:119
$insertRange = $ws.Range("A53:B61")
$insertRange.Insert(-4121)  # xlShiftDown

# Populate the newly inserted rows with the new content (English / Japanese pairs)
$ws.Range("A53").Value = 'Happy New Year'
$ws.Range("B53").Value = 'あけましておめでとうございます'
$ws.Range("A54").Value = 'Happy New Year (kanji)'
$ws.Range("B54").Value = '謹賀新年|きんがしんねん'
$ws.Range("A55").Value = 'Thank you for all your kind help during the past year.'
$ws.Range("B55").Value = '昨年は大変お世話になりました|さくねんはたいへんおせわになりました'
$ws.Range("A56").Value = 'I hope for your continued good will this year.'
$ws.Range("B56").Value = '本年もどうぞよろしくお願いいたします|ほんねんもどうぞよろしくおねがいいたします'
$ws.Range("A57").Value = 'I hope you are keeping well during the hot weather.'
$ws.Range("B57").Value = '暑中お見舞い申し上げます|しょちゅうおみまいもうしあげます'
$ws.Range("A58").Value = 'Congratulations on your graduation.'
$ws.Range("B58").Value = 'ご卒業おめでとうございます|ごそつぎょうおめでとうございます'
$ws.Range("A59").Value = 'Congratulations on your marriage.'
$ws.Range("B59").Value = 'ご結婚おめでとうございます|ごけっこんおめでとうございます'
$ws.Range("A60").Value = 'Happy Birthday'
$ws.Range("B60").Value = '誕生日おめでとう|たんじょうびおめでとう'
$ws.Range("A61").Value = 'Get well soon.'
$ws.Range("B61").Value = '早くよくなってください|はやくよくなってください'

